$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022252759893007
$ws.Range("D2").Value = 1.030752338702968
$ws.Range("E2").Value = 1.02302102075481
$ws.Range("F2").Value = 1.020697650145091
$ws.Range("I2").Value = 1.032912691336351
$ws.Range("J2").Value = 1.027439639688858
$ws.Range("K2").Value = 1.033562440400202
$ws.Range("L2").Value = 1.025853685629013
$ws.Range("M2").Value = 1.023537168260061
$ws.Range("N2").Value = 1.028898721826192

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023249111908376
$ws.Range("D3").Value = 1.031256055203984
$ws.Range("E3").Value = 1.023867409682729
$ws.Range("F3").Value = 1.022332814411987
$ws.Range("I3").Value = 1.033110015714587
$ws.Range("J3").Value = 1.028073717667712
$ws.Range("K3").Value = 1.033875395810438
$ws.Range("L3").Value = 1.026506756700713
$ws.Range("M3").Value = 1.024976355960243
$ws.Range("N3").Value = 1.029533700268504

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023893463459748
$ws.Range("D4").Value = 1.031581708288471
$ws.Range("E4").Value = 1.024415161479446
$ws.Range("F4").Value = 1.023390400685407
$ws.Range("I4").Value = 1.033236056739292
$ws.Range("J4").Value = 1.028483091010765
$ws.Range("K4").Value = 1.0340768979839
$ws.Range("L4").Value = 1.026928780453009
$ws.Range("M4").Value = 1.025906677666045
$ws.Range("N4").Value = 1.02994365496866

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024164264857275
$ws.Range("D5").Value = 1.031718543382869
$ws.Range("E5").Value = 1.024645456066863
$ws.Range("F5").Value = 1.02383490322629
$ws.Range("I5").Value = 1.033288651316737
$ws.Range("J5").Value = 1.028654972504056
$ws.Range("K5").Value = 1.034161369243622
$ws.Range("L5").Value = 1.027106066306928
$ws.Range("M5").Value = 1.026297568410138
$ws.Range("N5").Value = 1.030115780553386

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024209728705519
$ws.Range("D6").Value = 1.031741514489152
$ws.Range("E6").Value = 1.024684124704321
$ws.Range("F6").Value = 1.023909531098296
$ws.Range("I6").Value = 1.033297459122005
$ws.Range("J6").Value = 1.028683819334377
$ws.Range("K6").Value = 1.034175538232566
$ws.Range("L6").Value = 1.027135825597943
$ws.Range("M6").Value = 1.026363188161322
$ws.Range("N6").Value = 1.030144668349515

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023897082250147
$ws.Range("D7").Value = 1.031583536960164
$ws.Range("E7").Value = 1.02441823860923
$ws.Range("F7").Value = 1.023396340561016
$ws.Range("I7").Value = 1.033236761055369
$ws.Range("J7").Value = 1.028485388559732
$ws.Range("K7").Value = 1.03407802763786
$ws.Range("L7").Value = 1.026931149877131
$ws.Range("M7").Value = 1.025911901610701
$ws.Range("N7").Value = 1.029945955780409

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02258955576641
$ws.Range("D8").Value = 1.030922630387698
$ws.Range("E8").Value = 1.023307045094566
$ws.Range("F8").Value = 1.021250363899914
$ws.Range("I8").Value = 1.032979717704037
$ws.Range("J8").Value = 1.027654119646978
$ws.Range("K8").Value = 1.033668412010577
$ws.Range("L8").Value = 1.026074509604016
$ws.Range("M8").Value = 1.024023744386487
$ws.Range("N8").Value = 1.029113506370452

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020282766036584
$ws.Range("D9").Value = 1.029755896427656
$ws.Range("E9").Value = 1.021349594621962
$ws.Range("F9").Value = 1.017464939583013
$ws.Range("I9").Value = 1.032514214374268
$ws.Range("J9").Value = 1.026182263222738
$ws.Range("K9").Value = 1.032938978726784
$ws.Range("L9").Value = 1.024560716439517
$ws.Range("M9").Value = 1.020689190778867
$ws.Range("N9").Value = 1.027639559741265

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018742985900804
$ws.Range("D10").Value = 1.028976709015724
$ws.Range("E10").Value = 1.020045023161906
$ws.Range("F10").Value = 1.014938245282607
$ws.Range("I10").Value = 1.032195443763205
$ws.Range("J10").Value = 1.025196241730306
$ws.Range("K10").Value = 1.032447590098265
$ws.Range("L10").Value = 1.023548612598819
$ws.Range("M10").Value = 1.018460822289531
$ws.Range("N10").Value = 1.026652137985215

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018075773399514
$ws.Range("D11").Value = 1.028639002662249
$ws.Range("E11").Value = 1.019480218180555
$ws.Range("F11").Value = 1.013843332714447
$ws.Range("I11").Value = 1.032055413401451
$ws.Range("J11").Value = 1.024768138945929
$ws.Range("K11").Value = 1.032233610898374
$ws.Range("L11").Value = 1.023109663474649
$ws.Range("M11").Value = 1.017494566212241
$ws.Range("N11").Value = 1.02622342724578

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0178278677981
$ws.Range("D12").Value = 1.028513517549566
$ws.Range("E12").Value = 1.019270436409968
$ws.Range("F12").Value = 1.013436499318009
$ws.Range("I12").Value = 1.032003099374187
$ws.Range("J12").Value = 1.024608948939279
$ws.Range("K12").Value = 1.032153949024186
$ws.Range("L12").Value = 1.022946512125434
$ws.Range("M12").Value = 1.017135444459043
$ws.Range("N12").Value = 1.026064011171063

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017881047725306
$ws.Range("D13").Value = 1.02854043659346
$ws.Range("E13").Value = 1.019315434779937
$ws.Range("F13").Value = 1.013523772665904
$ws.Range("I13").Value = 1.032014334508889
$ws.Range("J13").Value = 1.024643103605964
$ws.Range("K13").Value = 1.032171044925142
$ws.Range("L13").Value = 1.022981513456471
$ws.Range("M13").Value = 1.017212486968956
$ws.Range("N13").Value = 1.026098214341291

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018055282946123
$ws.Range("D14").Value = 1.028628630954488
$ws.Range("E14").Value = 1.019462877313365
$ws.Range("F14").Value = 1.013809706523289
$ws.Range("I14").Value = 1.032051095237816
$ws.Range("J14").Value = 1.024754983796669
$ws.Range("K14").Value = 1.032227029701617
$ws.Range("L14").Value = 1.023096179501795
$ws.Range("M14").Value = 1.017464885446942
$ws.Range("N14").Value = 1.0262102534147

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.01816262527744
$ws.Range("D15").Value = 1.028682964356633
$ws.Range("E15").Value = 1.019553723091399
$ws.Range("F15").Value = 1.013985861787988
$ws.Range("I15").Value = 1.03207370491622
$ws.Range("J15").Value = 1.024823893846123
$ws.Range("K15").Value = 1.032261499871172
$ws.Range("L15").Value = 1.023166814955994
$ws.Range("M15").Value = 1.017620368273356
$ws.Range("N15").Value = 1.026279261324328

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01878725641448
$ws.Range("D16").Value = 1.028999114972329
$ws.Range("E16").Value = 1.020082509113399
$ws.Range("F16").Value = 1.015010892673793
$ws.Range("I16").Value = 1.032204694971479
$ws.Range("J16").Value = 1.025224629213771
$ws.Range("K16").Value = 1.032461765856215
$ws.Range("L16").Value = 1.023577729366756
$ws.Range("M16").Value = 1.018524920299543
$ws.Range("N16").Value = 1.026680565782163

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019178942171207
$ws.Range("D17").Value = 1.029197344864533
$ws.Range("E17").Value = 1.020414224483366
$ws.Range("F17").Value = 1.01565363744601
$ws.Range("I17").Value = 1.032286325996293
$ws.Range("J17").Value = 1.025475691643233
$ws.Range("K17").Value = 1.032587065259081
$ws.Range("L17").Value = 1.023835296816339
$ws.Range("M17").Value = 1.019091953274534
$ws.Range("N17").Value = 1.026931984749072

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019407359866868
$ws.Range("D18").Value = 1.029312938728078
$ws.Range("E18").Value = 1.020607716498035
$ws.Range("F18").Value = 1.016028459374702
$ws.Range("I18").Value = 1.032333746961765
$ws.Range("J18").Value = 1.025622021316132
$ws.Range("K18").Value = 1.03266003401244
$ws.Range("L18").Value = 1.023985463917455
$ws.Range("M18").Value = 1.019422563492724
$ws.Range("N18").Value = 1.027078522226891

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01948523664249
$ws.Range("D19").Value = 1.029352348078929
$ws.Range("E19").Value = 1.020673693656433
$ws.Range("F19").Value = 1.016156250611638
$ws.Range("I19").Value = 1.032349883542656
$ws.Range("J19").Value = 1.025671897189158
$ws.Range("K19").Value = 1.03268489473461
$ws.Range("L19").Value = 1.024036655567707
$ws.Range("M19").Value = 1.019535271113728
$ws.Range("N19").Value = 1.027128468929377

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019136922762282
$ws.Range("D20").Value = 1.029176079797006
$ws.Range("E20").Value = 1.020378633734331
$ws.Range("F20").Value = 1.015584685371856
$ws.Range("I20").Value = 1.032277587725463
$ws.Range("J20").Value = 1.025448766487244
$ws.Range("K20").Value = 1.03257363384135
$ws.Range("L20").Value = 1.023807669267291
$ws.Range("M20").Value = 1.019031129554617
$ws.Range("N20").Value = 1.026905021356273

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018003977047679
$ws.Range("D21").Value = 1.028602661172684
$ws.Range("E21").Value = 1.01941945883672
$ws.Range("F21").Value = 1.013725509898428
$ws.Range("I21").Value = 1.032040278414738
$ws.Range("J21").Value = 1.024722042668167
$ws.Range("K21").Value = 1.032210548558129
$ws.Range("L21").Value = 1.023062416144078
$ws.Range("M21").Value = 1.017390566275091
$ws.Range("N21").Value = 1.026177265506017

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017291225051058
$ws.Range("D22").Value = 1.028241864882576
$ws.Range("E22").Value = 1.018816456359184
$ws.Range("F22").Value = 1.012555792366041
$ws.Range("I22").Value = 1.031889333958357
$ws.Range("J22").Value = 1.024264117849401
$ws.Range("K22").Value = 1.031981218524016
$ws.Range("L22").Value = 1.022593231702202
$ws.Range("M22").Value = 1.016357853344496
$ws.Range("N22").Value = 1.025718690381484

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017669109045916
$ws.Range("D23").Value = 1.028433154601603
$ws.Range("E23").Value = 1.019136113037261
$ws.Range("F23").Value = 1.013175958502764
$ws.Range("I23").Value = 1.031969517262338
$ws.Range("J23").Value = 1.024506968016171
$ws.Range("K23").Value = 1.032102889494101
$ws.Range("L23").Value = 1.022842013714107
$ws.Range("M23").Value = 1.016905432614176
$ws.Range("N23").Value = 1.025961885423346

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019155909666095
$ws.Range("D24").Value = 1.029185688659312
$ws.Range("E24").Value = 1.020394715636701
$ws.Range("F24").Value = 1.0156158420921
$ws.Range("I24").Value = 1.032281536770271
$ws.Range("J24").Value = 1.025460933148266
$ws.Range("K24").Value = 1.032579703280122
$ws.Range("L24").Value = 1.023820153176198
$ws.Range("M24").Value = 1.019058613577108
$ws.Range("N24").Value = 1.026917205295349

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020879460209492
$ws.Range("D25").Value = 1.030057770656934
$ws.Range("E25").Value = 1.02185557143502
$ws.Range("F25").Value = 1.018444073422341
$ws.Range("I25").Value = 1.032636045560172
$ws.Range("J25").Value = 1.026563613537125
$ws.Range("K25").Value = 1.033128456164862
$ws.Range("L25").Value = 1.024952578759847
$ws.Range("M25").Value = 1.021552165416754
$ws.Range("N25").Value = 1.02802145161684

